$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.328.73'
$ws.Range("E2").Value = '  +1.82%  '

$ws.Range("D3").Value = '2.664.57'
$ws.Range("E3").Value = '  +1.89%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.17'
$ws.Range("E5").Value = '  -0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.45'
$ws.Range("E6").Value = '  -0.84%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("E8").Value = '  -0.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.173'
$ws.Range("E9").Value = '  +5.24%  '

$ws.Range("D10").Value = '2.663.83'
$ws.Range("E10").Value = '  +1.91%  '

$ws.Range("E11").Value = '  +2.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.357'
$ws.Range("E12").Value = '  +3.12%  '

$ws.Range("E13").Value = '  -0.08%  '

$ws.Range("D14").Value = '3.153.37'
$ws.Range("E14").Value = '  +1.80%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000189'
$ws.Range("E15").Value = '  +3.45%  '

$ws.Range("D16").Value = '72.258.09'
$ws.Range("E16").Value = '  +1.83%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.62'
$ws.Range("E17").Value = '  +0.17%  '

$ws.Range("D18").Value = '2.669.97'
$ws.Range("E18").Value = '  +1.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.99'
$ws.Range("E19").Value = '  +4.63%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.04'
$ws.Range("E20").Value = '  +3.86%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '379.01'
$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("E22").Value = '  +1.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.07'
$ws.Range("E23").Value = '  +11.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.58'
$ws.Range("E24").Value = '  +0.80%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.40'
$ws.Range("E25").Value = '  -0.59%  '

$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("E27").Value = '  +4.32%  '

$ws.Range("D28").Value = '2.807.13'
$ws.Range("E28").Value = '  +2.68%  '

$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("D30").Value = '0.0₃0950'
$ws.Range("E30").Value = '  -0.23%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.23'
$ws.Range("E31").Value = '  +2.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '522.76'
$ws.Range("E32").Value = '  -1.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.31'
$ws.Range("E33").Value = '  -0.36%  '

$ws.Range("E34").Value = '  -0.20%  '

$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.38'
$ws.Range("E36").Value = '  -1.26%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.55'
$ws.Range("E37").Value = '  +2.16%  '

$ws.Range("E38").Value = '  +0.80%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.39'
$ws.Range("E39").Value = '  +1.87%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.111'
$ws.Range("E40").Value = '  -5.98%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.84'
$ws.Range("E41").Value = '  -2.06%  '

$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.07'
$ws.Range("E43").Value = '  +1.23%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.60'
$ws.Range("E44").Value = '  -0.61%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.335'
$ws.Range("E45").Value = '  +1.28%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.33'
$ws.Range("E46").Value = '  -1.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '152.98'
$ws.Range("E47").Value = '  -0.63%  '

$ws.Range("E48").Value = '  +2.84%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.551'
$ws.Range("E49").Value = '  +4.03%  '

$ws.Range("E50").Value = '  +3.06%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0764'
$ws.Range("E51").Value = '  +1.66%  '
